$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.908.34'
$ws.Range("D3").Value = '1.636.32'
$ws.Range("E3").Value = '  -1.66%  '
$ws.Range("D4").Value = "'1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'214.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.80%  '
$ws.Range("D6").Value = "'0.5031"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.47%  '
$ws.Range("D7").Value = "'1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.38%  '
$ws.Range("D8").Value = "'0.06440"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").Value = "'0.2570"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("D10").Value = "'19.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.03%  '
$ws.Range("D11").Value = "'0.07722"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.84%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.643.81'
$ws.Range("E12").Value = '  -1.34%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'4.243"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.29%  '
$ws.Range("D14").Value = '1.863.97'
$ws.Range("E14").Value = '  -1.62%  '
$ws.Range("D15").Value = "'0.5444"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.69%  '
$ws.Range("D16").Value = '0.0₅7945'
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").Value = "'63.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.52%  '
$ws.Range("D18").Value = '25.929.61'
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("D19").Value = "'1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.20%  '
$ws.Range("D20").Value = "'204.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.71%  '
$ws.Range("D21").Value = "'4.301"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.97%  '
$ws.Range("D22").Value = "'9.966"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.36%  '
$ws.Range("D23").Value = "'5.946"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.98%  '
$ws.Range("D24").Value = "'1.007"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("D25").Value = "'1.894"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.92%  '
$ws.Range("D26").Value = "'141.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.01%  '
$ws.Range("D27").Value = "'0.1154"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.48%  '
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("D29").Value = "'6.730"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.28%  '
$ws.Range("D30").Value = "'0.05054"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.76%  '
$ws.Range("D31").Value = "'1.236"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.04%  '
$ws.Range("D32").Value = "'3.257"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'3.183"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.24%  '
$ws.Range("D34").Value = "'1.536"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.75%  '
$ws.Range("D35").Value = "'2.334"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.91%  '
$ws.Range("D36").Value = "'0.8925"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.13%  '
$ws.Range("D37").Value = "'2.618"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.96%  '
$ws.Range("D38").Value = "'0.5648"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.34%  '
$ws.Range("D39").Value = '1.152.02'
$ws.Range("E39").Value = '  -0.09%  '
$ws.Range("D40").Value = "'0.01565"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.85%  '
$ws.Range("D41").Value = "'2.554"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.18%  '
$ws.Range("D42").Value = "'1.004"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.45%  '
$ws.Range("D43").Value = "'5.631"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.8156"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.98%  '
$ws.Range("D45").Value = "'99.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.37%  '
$ws.Range("D46").Value = '1.778.42'
$ws.Range("E46").Value = '  -1.47%  '
$ws.Range("D47").Value = '0.0₈113'
$ws.Range("E47").Value = '  +2.54%  '
$ws.Range("D48").Value = "'0.4512"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.16%  '
$ws.Range("D49").Value = "'1.006"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").Value = "'54.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.34%  '
$ws.Range("D51").Value = "'0.05031"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.34%  '
